$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a question mark to the existing question in B43
$ws.Range("B43").Value = "ข้อ 2 ตอนเที่ยงวันพระจันทร์ไปไหน?"

# Add new row 44 with the new tag/response pair
$ws.Range("A44").Value = "ถูกสอง"
$ws.Range("B44").Value = "ข้อ 3 ปีอะไรเอ่ย มีหลากสี?"

# Update selection to match the new active cell
$ws.Range("B44").Select()
